# Add "Status" column to the student sheet, populate additional students,
# and refresh statuses/ages to match the new roster.
# (commit: "success : add student and get student endpoint")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add the new "Status" column ---
$ws.Cells.Item(1, 4).Value = "Status"

# --- Status column (column D), written first so the new shared strings
#     land in first-occurrence order: aktif, nonaktif, lulus, alumni ---
$statuses = @(
    "aktif", "aktif", "aktif", "aktif", "aktif",
    "nonaktif", "nonaktif", "nonaktif",
    "lulus", "lulus",
    "alumni", "alumni"
)
for ($i = 0; $i -lt $statuses.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $statuses[$i]
}

# --- Name column (column B), written next so the student names land in
#     the expected shared-string order right after the status words ---
$names = @(
    "Emma Smith", "Liam Johnson", "Olivia Williams", "Noa Brown",
    "Aiden Jones", "Amelia Miller", "Ethan Davis", "Charlotte Garcia",
    "Harper Rodriguez", "Mason Wilson", "Scarlett Moore", "Logan Taylor"
)
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $names[$i]
}

# --- Age column (column C) ---
$ages = @(17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28)
for ($i = 0; $i -lt $ages.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $ages[$i]
}

# --- ID column (column A) ---
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}

# --- Column B width (auto-fit to the longest name) ---
$ws.Columns.Item(2).AutoFit()

# --- Selection / active cell lands on B14, just past the last data row ---
[void]$ws.Range("B14").Select()

# --- Drop the explicit print scaling so the page just prints at 100% ---
$ws.PageSetup.Zoom = $true

Write-Output "done"
